$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 3744
$ws.Range("F5").Value = 3744
$ws.Range("F7").Value = 5280
$ws.Range("F8").Value = 582
$ws.Range("F9").Value = 414
$ws.Range("F10").Value = 222
$ws.Range("F11").Value = 1047
$ws.Range("F13").Value = 134
$ws.Range("F14").Value = 47
$ws.Range("F15").Value = 726
$ws.Range("F16").Value = 356
$ws.Range("F19").Value = 170
$ws.Range("F22").Value = 6039
$ws.Range("F23").Value = 6039
$ws.Range("F27").Value = 6798
$ws.Range("F28").Value = 23
$ws.Range("F32").Value = 747
$ws.Range("F34").Value = 323
$ws.Range("F36").Value = 150
$ws.Range("F37").Value = 1133
$ws.Range("F41").Value = 918
$ws.Range("F42").Value = 1116

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 1150

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 1150
$ws.Range("F7").Value = 3744
$ws.Range("F8").Value = 3744
$ws.Range("F10").Value = 5280
$ws.Range("F11").Value = 582
$ws.Range("F12").Value = 414
$ws.Range("F13").Value = 222
$ws.Range("F14").Value = 1047
$ws.Range("F16").Value = 134
$ws.Range("F17").Value = 47
$ws.Range("F18").Value = 726
$ws.Range("F19").Value = 356
$ws.Range("F23").Value = 170
$ws.Range("F26").Value = 6039
$ws.Range("F30").Value = 6798
$ws.Range("F31").Value = 23
$ws.Range("F35").Value = 747
$ws.Range("F37").Value = 323
$ws.Range("F40").Value = 150
$ws.Range("F41").Value = 1133
$ws.Range("F45").Value = 918
$ws.Range("F46").Value = 1116
